$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the bugged values used to compute the effective wind speed ---
$ws.Range("B6").Value = 18
$ws.Range("C6").Value = 45

# --- Remove the old duplicate "Equation 6" computation living in column G / I ---
$ws.Range("G3:G4").Clear()
$ws.Range("I6").Clear()

# --- Row 9 now holds the (fixed) Uws calculation, with a bold-ish header style ---
$ws.Range("A9").Value = "Uws"
$ws.Range("A9").Font.Name = "Calibri"

# Fix the suppression-map bug: SIN(C6^2) -> SIN(C6)^2
$ws.Range("B9").Formula = "=+B6*(((A6/B6)^2+2*(A6/B6)*SIN(C6)*COS(D6)+SIN(C6)^2)^0.5)"

# --- Move the "EffectiveWindSpeed" label down to row 10, column A ---
$ws.Range("A10").Value = "EffectiveWindSpeed"

# --- Restore the active cell/selection shown when the file was last saved ---
$ws.Range("F6").Select()
